# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# The old scraping code only pulled team statistics, not the season record,
# so this backfills the three new columns (AD, AE, AF) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column titles.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold font, border, centered/top
# aligned) by copying the format only from an existing header cell (A1)
# onto the new header cells, leaving the values we just set intact.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-54): every player on the roster shares the same team season
# record: 85 wins, 77 losses, 0 ties.
$ws.Range("AD2:AD54").Value = 85
$ws.Range("AE2:AE54").Value = 77
$ws.Range("AF2:AF54").Value = 0
